$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 "Marking": Right count and Wrong penalty corrections
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": corrected totals and the Max summary text
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "50 / 112"
